$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$newValues = @{
    2  = 0.95
    3  = 1.19
    4  = 0.86
    5  = 0.98
    6  = 0.83
    7  = 0.92
    8  = 0.87
    9  = 1.17
    10 = 0.9
    11 = 0.82
    12 = 0.85
    13 = 1.02
    14 = 0.84
    15 = 1.02
    16 = 0.88
    17 = 1.07
}

foreach ($row in $newValues.Keys) {
    $ws.Range("L$row").Value = $newValues[$row]
}
